# Add new columns I (I0) and J (IF) to the sheet, matching the H (IP) column's
# header style, then fill in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the formatting from H1 (bold, bordered,
#     centered header style) onto I1/J1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-7): plain numeric values, no special formatting
#     (matches the un-styled data cells in the other columns).
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 5

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 6
